$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 6232382
$ws.Range("I100").Value = 11149372
$ws.Range("J100").Value = 86144.164
$ws.Range("K100").Value = 11149372
$ws.Range("L100").Value = 86144.164
$ws.Range("M100").Value = -11148831
$ws.Range("N100").Value = -87226.164

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 31833.334
$ws.Range("J136").Value = 31833.334
$ws.Range("L136").Value = 31833.334
$ws.Range("N136").Value = -42033.334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 76924170
$ws.Range("I137").Value = 83334440
$ws.Range("J137").Value = 900
$ws.Range("K137").Value = 250003320
$ws.Range("L137").Value = 2700
$ws.Range("M137").Value = -250000770
$ws.Range("N137").Value = -7800

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1639.9412
$ws.Range("I2").Value = 646.2727
$ws.Range("J2").Value = 3461.6667
$ws.Range("K2").Value = 646.2727
$ws.Range("L2").Value = 3461.6667
$ws.Range("M2").Value = -533.2727
$ws.Range("N2").Value = -3687.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1639.9412
$ws.Range("I116").Value = 646.2727
$ws.Range("J116").Value = 3461.6667
$ws.Range("K116").Value = 646.2727
$ws.Range("L116").Value = 3461.6667
$ws.Range("M116").Value = 1647.7273
$ws.Range("N116").Value = -8049.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3015.111
$ws.Range("I122").Value = 3403.8333
$ws.Range("K122").Value = 10211.4999
$ws.Range("M122").Value = -7761.499899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value = 63000
$ws.Range("J138").Value = 63000
$ws.Range("L138").Value = 63000
$ws.Range("N138").Value = -73280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1639.9412
$ws.Range("I3").Value = 646.2727
$ws.Range("J3").Value = 3461.6667
$ws.Range("K3").Value = 646.2727
$ws.Range("L3").Value = 3461.6667
$ws.Range("M3").Value = -532.2727
$ws.Range("N3").Value = -3689.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 615.3077
$ws.Range("I80").Value = 238
$ws.Range("J80").Value = 664.5217
$ws.Range("K80").Value = 238
$ws.Range("L80").Value = 664.5217
$ws.Range("M80").Value = 760
$ws.Range("N80").Value = -2660.5217

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 615.3077
$ws.Range("I83").Value = 238
$ws.Range("J83").Value = 664.5217
$ws.Range("K83").Value = 1190
$ws.Range("L83").Value = 3322.6085
$ws.Range("M83").Value = 3802
$ws.Range("N83").Value = -13306.6085

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 380.46155
$ws.Range("I107").Value = 216.83333
$ws.Range("J107").Value = 520.7143
$ws.Range("K107").Value = 216.83333
$ws.Range("L107").Value = 520.7143
$ws.Range("M107").Value = 1703.16667
$ws.Range("N107").Value = -4360.7143

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 14577.5
$ws.Range("J122").Value = 14577.5
$ws.Range("L122").Value = 14577.5
$ws.Range("N122").Value = -24377.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 45935.363
$ws.Range("I16").Value = 62826.25
$ws.Range("J16").Value = 893
$ws.Range("K16").Value = 62826.25
$ws.Range("L16").Value = 893
$ws.Range("M16").Value = -62539.25
$ws.Range("N16").Value = -1467

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1154.7142
$ws.Range("I31").Value = 1089.6923
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 1089.6923
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = -794.6922999999999
$ws.Range("N31").Value = -2590

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1154.7142
$ws.Range("I34").Value = 1089.6923
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 1089.6923
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = -887.6922999999999
$ws.Range("N34").Value = -2404

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H110").Value = 50000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 50000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 50000
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -58180

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 45935.363
$ws.Range("I113").Value = 62826.25
$ws.Range("J113").Value = 893
$ws.Range("K113").Value = 62826.25
$ws.Range("L113").Value = 893
$ws.Range("M113").Value = -60656.25
$ws.Range("N113").Value = -5233

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1918
$ws.Range("I122").Value = 1334.1666
$ws.Range("J122").Value = 2554.9092
$ws.Range("K122").Value = 4002.4998
$ws.Range("L122").Value = 7664.7276
$ws.Range("M122").Value = -1552.4998
$ws.Range("N122").Value = -12564.7276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3857.3333
$ws.Range("I132").Value = 2541.5
$ws.Range("J132").Value = 5361.143
$ws.Range("K132").Value = 7624.5
$ws.Range("L132").Value = 16083.429
$ws.Range("M132").Value = -5094.5
$ws.Range("N132").Value = -21143.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1016.1053
$ws.Range("I5").Value = 1016.1053
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 3048.3159
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -2936.3159
$ws.Range("N5").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4608.25
$ws.Range("I56").Value = 4608.25
$ws.Range("K56").Value = 4608.25
$ws.Range("M56").Value = -4078.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1671.0769
$ws.Range("J131").Value = 1871.2727
$ws.Range("L131").Value = 5613.8181
$ws.Range("N131").Value = -15693.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1441.2142
$ws.Range("J132").Value = 1475.1538
$ws.Range("L132").Value = 13276.3842
$ws.Range("N132").Value = -18336.3842

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1016.1053
$ws.Range("I135").Value = 1016.1053
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 9144.947700000001
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -6609.947700000001
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2262
$ws.Range("I113").Value = 2252.5
$ws.Range("J113").Value = 2300
$ws.Range("K113").Value = 2252.5
$ws.Range("L113").Value = 2300
$ws.Range("M113").Value = -82.5
$ws.Range("N113").Value = -6640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2840.4
$ws.Range("I132").Value = 1604.5714
$ws.Range("J132").Value = 5724
$ws.Range("K132").Value = 4813.7142
$ws.Range("L132").Value = 17172
$ws.Range("M132").Value = -2283.7142
$ws.Range("N132").Value = -22232

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7627.2666
$ws.Range("I22").Value = 810.625
$ws.Range("J22").Value = 15417.714
$ws.Range("K22").Value = 810.625
$ws.Range("L22").Value = 15417.714
$ws.Range("M22").Value = -515.625
$ws.Range("N22").Value = -16007.714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 7627.2666
$ws.Range("I27").Value = 810.625
$ws.Range("J27").Value = 15417.714
$ws.Range("K27").Value = 810.625
$ws.Range("L27").Value = 15417.714
$ws.Range("M27").Value = -703.625
$ws.Range("N27").Value = -15631.714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 37499.75
$ws.Range("J70").Value = 37499.75
$ws.Range("L70").Value = 37499.75
$ws.Range("N70").Value = -38129.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 37499.75
$ws.Range("J73").Value = 37499.75
$ws.Range("L73").Value = 37499.75
$ws.Range("N73").Value = -39683.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3920.0435
$ws.Range("J81").Value = 4911.7646
$ws.Range("L81").Value = 9823.529200000001
$ws.Range("N81").Value = -11945.5292

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 48749.75
$ws.Range("J82").Value = 48749.75
$ws.Range("L82").Value = 48749.75
$ws.Range("N82").Value = -49515.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3920.0435
$ws.Range("J84").Value = 4911.7646
$ws.Range("L84").Value = 49117.64600000001
$ws.Range("N84").Value = -59725.64600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H85").Value = 48749.75
$ws.Range("J85").Value = 48749.75
$ws.Range("L85").Value = 48749.75
$ws.Range("N85").Value = -51401.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13516544
$ws.Range("I132").Value = 18520714
$ws.Range("K132").Value = 55562142
$ws.Range("M132").Value = -55559612

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 285000
$ws.Range("J135").Value = 285000
$ws.Range("L135").Value = 285000
$ws.Range("N135").Value = -295140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 13375204
$ws.Range("I136").Value = 27862020
$ws.Range("J136").Value = 2758.1538
$ws.Range("K136").Value = 83586060
$ws.Range("L136").Value = 8274.4614
$ws.Range("M136").Value = -83583510
$ws.Range("N136").Value = -13374.4614
